$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = -0
$ws.Range("B2").Value = -0.07498804478836586
$ws.Range("C2").Value = -0
$ws.Range("D2").Value = 0.2123108505157703
$ws.Range("E2").Value = 0.005516432524130028
$ws.Range("G2").Value = 0
$ws.Range("I2").Value = -0
$ws.Range("J2").Value = -0
$ws.Range("K2").Value = 0.0210650617049058
$ws.Range("L2").Value = -0
$ws.Range("M2").Value = 0.2062088150060369
$ws.Range("N2").Value = -0.007387703804133994
$ws.Range("R2").Value = -0
$ws.Range("S2").Value = 0
$ws.Range("T2").Value = -0.09179633776493132
$ws.Range("V2").Value = 0.01626758625129003
$ws.Range("W2").Value = -0.03827681932504355
$ws.Range("Y2").Value = -0
$ws.Range("Z2").Value = -0
$ws.Range("AB2").Value = 0
$ws.Range("AC2").Value = -0.05304819134633224
$ws.Range("AD2").Value = 0
$ws.Range("AE2").Value = -0.01596718197831162
$ws.Range("AF2").Value = 0.001357506954384982
$ws.Range("AG2").Value = -0
$ws.Range("AH2").Value = -0
$ws.Range("AI2").Value = -0
$ws.Range("AJ2").Value = 0
$ws.Range("AK2").Value = -0
$ws.Range("AL2").Value = -0.03337734227611048
$ws.Range("AM2").Value = 0
$ws.Range("AN2").Value = 0.02710235950979624
$ws.Range("AO2").Value = 0.06807114299100285
$ws.Range("AQ2").Value = 0
$ws.Range("AR2").Value = -0
$ws.Range("AT2").Value = 0
$ws.Range("AU2").Value = -0.1493238582359526
$ws.Range("AW2").Value = 0.07989574898010814
$ws.Range("AX2").Value = 0.004724416540976223
$ws.Range("AY2").Value = -0
$ws.Range("BC2").Value = -0
$ws.Range("BD2").Value = -0.01731206187656933
$ws.Range("BF2").Value = 0.08423891762249715
$ws.Range("BG2").Value = 0.03181235655624312
$ws.Range("BJ2").Value = -0
$ws.Range("BL2").Value = 0
$ws.Range("BM2").Value = 0.03624722305829137
$ws.Range("BO2").Value = -0.03913293478256753
$ws.Range("BP2").Value = -0.08637971709137597
$ws.Range("BU2").Value = 0
$ws.Range("BV2").Value = -0.04500391369559421
$ws.Range("BW2").Value = 0
$ws.Range("BX2").Value = 0.01500758071599884
$ws.Range("BY2").Value = -0.02086773070508119
$ws.Range("BZ2").Value = -0
$ws.Range("CB2").Value = 0
$ws.Range("CD2").Value = -0
$ws.Range("CE2").Value = 0.03220467152992416
$ws.Range("CG2").Value = -0.03364726260395302
$ws.Range("CH2").Value = 0.01645725400885414
$ws.Range("CJ2").Value = -0
$ws.Range("CM2").Value = -0
$ws.Range("CN2").Value = -0.01104918755697818
$ws.Range("CP2").Value = 0.02074356515685865
$ws.Range("CQ2").Value = 0.03582172844414858
$ws.Range("CT2").Value = 0
$ws.Range("CU2").Value = -0
$ws.Range("CV2").Value = -0
$ws.Range("CW2").Value = 0.04603546316646197
$ws.Range("CY2").Value = -0.03621017246313984
$ws.Range("CZ2").Value = 0.009757611166591691
$ws.Range("DD2").Value = -0
$ws.Range("DE2").Value = -0
$ws.Range("DF2").Value = 0.02952187127600957
$ws.Range("DH2").Value = 0.02516906389234852
$ws.Range("DI2").Value = 0.03783215140409813
$ws.Range("DJ2").Value = 0
$ws.Range("DK2").Value = -0
$ws.Range("DL2").Value = -0
$ws.Range("DN2").Value = 0
$ws.Range("DO2").Value = -0.01844302020456456
$ws.Range("DQ2").Value = 0.03409103206182655
$ws.Range("DR2").Value = -0.01622991875651391
$ws.Range("DS2").Value = -0
$ws.Range("DW2").Value = 0
$ws.Range("DX2").Value = -0.0547282908726129
$ws.Range("DY2").Value = -0
$ws.Range("DZ2").Value = -0.01156697219937911
$ws.Range("EA2").Value = -0.02323154868926323
$ws.Range("EB2").Value = 0
$ws.Range("EF2").Value = -0
$ws.Range("EG2").Value = 0.04188967657062639
$ws.Range("EI2").Value = 0.07034029962704759
$ws.Range("EJ2").Value = -0.02720351628130804
$ws.Range("EO2").Value = 0
$ws.Range("EP2").Value = 0.04467283160219787
$ws.Range("EQ2").Value = 0
$ws.Range("ER2").Value = -0.0332712947545084
$ws.Range("ES2").Value = 0.03374213308970007
$ws.Range("ET2").Value = 0
$ws.Range("EU2").Value = -0
$ws.Range("EV2").Value = 0
$ws.Range("EX2").Value = 0
$ws.Range("EY2").Value = 0.04336569382910841
$ws.Range("FA2").Value = -0.02930277627662451
$ws.Range("FB2").Value = 0.01670699430381313
$ws.Range("FD2").Value = -0
$ws.Range("FG2").Value = -0
$ws.Range("FH2").Value = 0.003066748497303185
$ws.Range("FI2").Value = 0
$ws.Range("FJ2").Value = -0.006813118420519349
$ws.Range("FK2").Value = -0.006099503072460796
$ws.Range("FL2").Value = -0
$ws.Range("FN2").Value = -0
$ws.Range("FP2").Value = -0
$ws.Range("FQ2").Value = -0.01240537822611858
$ws.Range("FR2").Value = -0
$ws.Range("FS2").Value = -0.02042875914825069
$ws.Range("FT2").Value = 0.007039864533620131
$ws.Range("FV2").Value = -0
$ws.Range("FW2").Value = -0
$ws.Range("FY2").Value = 0
$ws.Range("FZ2").Value = -0.03029568608595597
$ws.Range("GB2").Value = 0.03483274327951914
$ws.Range("GD2").Value = 0
$ws.Range("GE2").Value = -0
